$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 107, shifting the existing rows 107:133 down to 108:134.
$ws.Rows.Item(107).Insert()

# Populate the newly inserted row 107 with the new weekly data point.
$ws.Range("A107").Value = 11
$ws.Range("B107").Value = "Vega Monumental Concepción"
$ws.Range("C107").Value = "Bíobío"
$ws.Range("D107").Value = 44722
$ws.Range("E107").Value = 8
$ws.Range("F107").Value = 100112043
$ws.Range("G107").Value = "Pepino ensalada"
$ws.Range("H107").Value = "Sin especificar"
$ws.Range("I107").Value = "Primera"
$ws.Range("J107").Value = 300
$ws.Range("K107").Value = 19000
$ws.Range("L107").Value = 20000
$ws.Range("M107").Value = 19500
$ws.Range("N107").Value = "$/caja 60 unidades"
$ws.Range("O107").Value = "Región de Arica y Parinacota"
$ws.Range("P107").Value = 325
$ws.Range("Q107").Value = 60
$ws.Range("R107").Value = "Hortaliza"
